$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E single-cell swaps (missing <-> restored values) ---
$ws.Range("E5").Value = $null
$ws.Range("E8").Value = -6.6
$ws.Range("E12").Value = $null
$ws.Range("E14").Value = -5.4
$ws.Range("E18").Value = $null

# --- Remove row 26 ("RM 232") entirely; rows below shift up ---
$ws.Range("A26:F26").EntireRow.Delete()

# --- Remove the (now shifted) "SC 92" row entirely; rows below shift up again ---
$ws.Range("A27:F27").EntireRow.Delete()

# --- Fix up values within the remaining SC rows ---
$ws.Range("B26").Value = -20.2
$ws.Range("B27").Value = $null
$ws.Range("C33").Value = 10.4
